$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header (copy formatting from B1, which already has the
# bold / bordered / centered style used for all header & label cells)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

# Row 2..23: label (col A) + two value columns (B, C), reordered and
# updated according to the new export.
$rows = @(
    @("1810 Renta variable",                57522.86,             57731.17),
    @("1822 Raices Valores Negociables",     78555.98,             77809.84),
    @("Alpha Acciones",                      15483.59,             15206.27),
    @("Alpha Mega",                          165348.66,            164426.67),
    @("Arpenta acciones",                    5762.1,               5757.28),
    @("Balanz",                              57383.74,             56743.94),
    @("FBA Acciones Argentinas",             20281.88,             20277.42),
    @("FBA Calificado",                      22444.09,             22708.68),
    @("Fima Acciones",                       71634.56,             70690.85000000001),
    @("Fima PB Acciones",                    91072.45,             90675.48),
    @("Goal Acciones Argentinas",            151.45,               154.56),
    @("Goal acciones plus",                  1009.85,              988.5599999999999),
    @("HF Acciones Argentinas",              8200.209999999999,    7943.86),
    @("HF Acciones Lideres",                 74332.42999999999,    74052),
    @("Lombardi",                            3625.45,              3676.3),
    @("MAF",                                 4310.19,              4363.18),
    @("Pellegrini Acciones",                 20052.71,             19832.28),
    @("Pionero Acciones",                    17731.57,             17820.95),
    @("Premier Renta Variable",              11202.22,             11152.86),
    @("Supefondo RV",                        231878.3,             231829.55),
    @("avg",                                 47899.21,             47692.08),
    @("total",                               957984.29,            953841.7)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
